$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated E2H values for the "Central" rows (rows 2-10)
$ws.Range("B2").Value = 0.8319027428269387
$ws.Range("B3").Value = 0.4953909906185326
$ws.Range("B5").Value = 0.8552075670435977
$ws.Range("B6").Value = 0.9402498030831551
$ws.Range("B7").Value = 0.2122967701559634
$ws.Range("B9").Value = -0.2949992232077283
$ws.Range("B10").Value = -0.2184503140238934

# Remove the "Decentral" rows (previously rows 12-20); data now ends at row 11
$ws.Range("A12:B20").ClearContents() | Out-Null
